$wb = $excel.ActiveWorkbook

# --- module_due_dates sheet: relabel header/module names, add Capstone rows ---
$ws = $wb.Worksheets.Item("module_due_dates")

$ws.Range("A1").Value = "Assignment"
$ws.Range("A2").Value = "Module 1: Good and bad visualizations"
$ws.Range("A3").Value = "Module 2: Coding Fundamentals"
$ws.Range("A4").Value = "Module 3: Data Exploration"
$ws.Range("A5").Value = "Module 4: Putting it together"

$ws.Range("A6").Value = "Capstone plan"
$ws.Range("B6").Value = [DateTime]"2023-11-07"

$ws.Range("A7").Value = "Capstone"
$ws.Range("B7").Value = [DateTime]"2023-12-05"

# --- switch the active tab from Schedule to module_due_dates ---
$ws.Activate()
$ws.Range("B8").Select()
